# Add filament tip analysis row to the "Tabelle2" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

$ws.Range("B5").Value = "BU6981_2U_L.csv"
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = "906 22705"
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "'0.35"
$ws.Range("H5").Value = "'0.7"

# G5/H5 were entered as quoted text (to match "0.35"/"0.7" stored as text,
# like the existing F3:H4 cells) - re-pasting the formatting from a sibling
# cell clears the quote-prefix style bit Excel adds automatically, without
# touching the underlying General number format / style index.
$ws.Range("H4").Copy() | Out-Null
$ws.Range("G5:H5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Columns.Item(4).AutoFit() | Out-Null

$ws.Range("H6").Select()
